# Commit: "updated tests: use start.bat for testing instead of Are.exe
# because start.bat is now without console window and does not need an
# installed .NET framework"
#
# The only textual change in the workbook is that every "ARE start file: .."
# instruction cell (column C, rows 7-21 of the "Test Scenarios" sheet) that
# used to reference "Are.exe" now references "start.bat" instead. The rest
# of the cell text (including the embedded newline before "ACS: ACS.exe")
# stays identical.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C, rows 7 through 21 all start with a line like
#   "ARE start file: Are.exe (start.sh - Linux)"
# Replace the executable name in-place for every one of them, keeping the
# rest of the (possibly multi-line) text untouched.
for ($row = 7; $row -le 21; $row++) {
    $cell = $ws.Cells.Item($row, 3)   # column C
    $text = $cell.Text
    if ($text -and $text.Contains("Are.exe")) {
        $cell.Value = $text.Replace("Are.exe", "start.bat")
    }
}

# Reflect the saved view state: the author had scrolled the sheet down and
# had cell C21 selected when the workbook was last saved.
$win = $excel.ActiveWindow
$win.ScrollRow = 18
$win.ScrollColumn = 1
$ws.Range("C21").Select()
